# Auto-generated edit script applying the Famfrit_Profits diff
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H9").Value = 9938.071
$ws_ALC.Range("I9").Value = 17423.334
$ws_ALC.Range("J9").Value = 4324.125
$ws_ALC.Range("K9").Value = 17423.334
$ws_ALC.Range("L9").Value = 4324.125
$ws_ALC.Range("M9").Value = -17254.334
$ws_ALC.Range("N9").Value = -4662.125
$ws_ALC.Range("H21").Value = 0
$ws_ALC.Range("J21").Value = 0
$ws_ALC.Range("L21").Value = 0
$ws_ALC.Range("N21").ClearContents()
$ws_ALC.Range("H23").Value = 0
$ws_ALC.Range("J23").Value = 0
$ws_ALC.Range("L23").Value = 0
$ws_ALC.Range("N23").ClearContents()
$ws_ALC.Range("H29").Value = 1000
$ws_ALC.Range("J29").Value = 1500
$ws_ALC.Range("L29").Value = 4500
$ws_ALC.Range("N29").Value = -5062
$ws_ALC.Range("H40").Value = 0
$ws_ALC.Range("J40").Value = 0
$ws_ALC.Range("L40").Value = 0
$ws_ALC.Range("N40").ClearContents()
$ws_ALC.Range("H107").Value = 2100
$ws_ALC.Range("I107").Value = 0
$ws_ALC.Range("J107").Value = 2100
$ws_ALC.Range("K107").Value = 0
$ws_ALC.Range("L107").Value = 2100
$ws_ALC.Range("M107").ClearContents()
$ws_ALC.Range("N107").Value = -5940
$ws_ALC.Range("H138").Value = 4667.1665
$ws_ALC.Range("I138").Value = 2052.8462
$ws_ALC.Range("J138").Value = 6144.826
$ws_ALC.Range("K138").Value = 6158.5386
$ws_ALC.Range("L138").Value = 18434.478
$ws_ALC.Range("M138").Value = -1018.5386
$ws_ALC.Range("N138").Value = -28714.478
$ws_ARM.Range("H133").Value = 92703.375
$ws_ARM.Range("J133").Value = 92703.375
$ws_ARM.Range("L133").Value = 92703.375
$ws_ARM.Range("N133").Value = -97763.375
$ws_ARM.Range("H134").Value = 79998.664
$ws_ARM.Range("J134").Value = 79998.664
$ws_ARM.Range("L134").Value = 79998.664
$ws_ARM.Range("N134").Value = -90138.664
$ws_CRP.Range("H58").Value = 1402
$ws_CRP.Range("I58").Value = 1416.6666
$ws_CRP.Range("K58").Value = 1416.6666
$ws_CRP.Range("M58").Value = -1213.6666
$ws_CRP.Range("H99").Value = 8748.9
$ws_CRP.Range("I99").Value = 8642.857
$ws_CRP.Range("K99").Value = 8642.857
$ws_CRP.Range("M99").Value = -7144.857
$ws_CRP.Range("H126").Value = 8748.9
$ws_CRP.Range("I126").Value = 8642.857
$ws_CRP.Range("K126").Value = 25928.571
$ws_CRP.Range("M126").Value = -23458.571
$ws_CRP.Range("H136").Value = 1402
$ws_CRP.Range("I136").Value = 1416.6666
$ws_CRP.Range("K136").Value = 4249.9998
$ws_CRP.Range("M136").Value = -1699.9998
$ws_CRP.Range("H140").Value = 133332
$ws_CRP.Range("J140").Value = 133332
$ws_CRP.Range("L140").Value = 133332
$ws_CRP.Range("N140").Value = -143692
$ws_CRP.Range("H141").Value = 93408.914
$ws_CRP.Range("J141").Value = 99264.27
$ws_CRP.Range("L141").Value = 99264.27
$ws_CRP.Range("N141").Value = -109624.27
$ws_CUL.Range("H5").Value = 3730.625
$ws_CUL.Range("J5").Value = 4152.5
$ws_CUL.Range("L5").Value = 12457.5
$ws_CUL.Range("N5").Value = -12681.5
$ws_CUL.Range("H92").Value = 1536.8182
$ws_CUL.Range("I92").Value = 1186
$ws_CUL.Range("J92").Value = 1957.8
$ws_CUL.Range("K92").Value = 3558
$ws_CUL.Range("L92").Value = 5873.4
$ws_CUL.Range("M92").Value = -2310
$ws_CUL.Range("N92").Value = -8369.4
$ws_CUL.Range("H109").Value = 5270.8184
$ws_CUL.Range("I109").Value = 13999.667
$ws_CUL.Range("J109").Value = 1997.5
$ws_CUL.Range("K109").Value = 41999.001
$ws_CUL.Range("L109").Value = 5992.5
$ws_CUL.Range("M109").Value = -40959.001
$ws_CUL.Range("N109").Value = -8072.5
$ws_CUL.Range("H116").Value = 2500
$ws_CUL.Range("I116").Value = 2500
$ws_CUL.Range("K116").Value = 7500
$ws_CUL.Range("M116").Value = -4058
$ws_CUL.Range("H121").Value = 364367.62
$ws_CUL.Range("I121").Value = 1336
$ws_CUL.Range("J121").Value = 571814.3
$ws_CUL.Range("K121").Value = 4008
$ws_CUL.Range("L121").Value = 1715442.9
$ws_CUL.Range("M121").Value = -2698
$ws_CUL.Range("N121").Value = -1718062.9
$ws_CUL.Range("H122").Value = 1428.4615
$ws_CUL.Range("I122").Value = 459
$ws_CUL.Range("J122").Value = 1719.3
$ws_CUL.Range("K122").Value = 4131
$ws_CUL.Range("L122").Value = 15473.7
$ws_CUL.Range("M122").Value = -1681
$ws_CUL.Range("N122").Value = -20373.7
$ws_CUL.Range("H131").Value = 25001240
$ws_CUL.Range("I131").Value = 41667516
$ws_CUL.Range("J131").Value = 1828.5
$ws_CUL.Range("K131").Value = 125002548
$ws_CUL.Range("L131").Value = 5485.5
$ws_CUL.Range("M131").Value = -124997508
$ws_CUL.Range("N131").Value = -15565.5
$ws_CUL.Range("H133").Value = 2610.8333
$ws_CUL.Range("I133").Value = 1566.3334
$ws_CUL.Range("K133").Value = 4699.0002
$ws_CUL.Range("M133").Value = 360.9997999999996
$ws_CUL.Range("H135").Value = 3730.625
$ws_CUL.Range("J135").Value = 4152.5
$ws_CUL.Range("L135").Value = 37372.5
$ws_CUL.Range("N135").Value = -42442.5
$ws_CUL.Range("H139").Value = 2624.9167
$ws_CUL.Range("I139").Value = 2489
$ws_CUL.Range("K139").Value = 7467
$ws_CUL.Range("M139").Value = -2327
$ws_GSM.Range("H80").Value = 4955
$ws_GSM.Range("I80").Value = 2360
$ws_GSM.Range("J80").Value = 7117.5
$ws_GSM.Range("K80").Value = 2360
$ws_GSM.Range("L80").Value = 7117.5
$ws_GSM.Range("M80").Value = -1362
$ws_GSM.Range("N80").Value = -9113.5
$ws_GSM.Range("H83").Value = 4955
$ws_GSM.Range("I83").Value = 2360
$ws_GSM.Range("J83").Value = 7117.5
$ws_GSM.Range("K83").Value = 11800
$ws_GSM.Range("L83").Value = 35587.5
$ws_GSM.Range("M83").Value = -6808
$ws_GSM.Range("N83").Value = -45571.5
$ws_LTW.Range("H7").Value = 4068.2856
$ws_LTW.Range("I7").Value = 3000
$ws_LTW.Range("J7").Value = 4869.5
$ws_LTW.Range("K7").Value = 3000
$ws_LTW.Range("L7").Value = 4869.5
$ws_LTW.Range("M7").Value = -2888
$ws_LTW.Range("N7").Value = -5093.5
$ws_LTW.Range("H20").Value = 6193.8823
$ws_LTW.Range("I20").Value = 6981.4546
$ws_LTW.Range("J20").Value = 4750
$ws_LTW.Range("K20").Value = 6981.4546
$ws_LTW.Range("L20").Value = 4750
$ws_LTW.Range("M20").Value = -6755.4546
$ws_LTW.Range("N20").Value = -5202
$ws_LTW.Range("H25").Value = 3253.75
$ws_LTW.Range("J25").Value = 4002.6667
$ws_LTW.Range("L25").Value = 4002.6667
$ws_LTW.Range("N25").Value = -4462.6667
$ws_LTW.Range("H46").Value = 2131.5833
$ws_LTW.Range("I46").Value = 989.96155
$ws_LTW.Range("J46").Value = 5099.8
$ws_LTW.Range("K46").Value = 989.96155
$ws_LTW.Range("L46").Value = 5099.8
$ws_LTW.Range("M46").Value = -801.96155
$ws_LTW.Range("N46").Value = -5475.8
$ws_LTW.Range("H126").Value = 4068.2856
$ws_LTW.Range("I126").Value = 3000
$ws_LTW.Range("J126").Value = 4869.5
$ws_LTW.Range("K126").Value = 9000
$ws_LTW.Range("L126").Value = 14608.5
$ws_LTW.Range("M126").Value = -6530
$ws_LTW.Range("N126").Value = -19548.5
$ws_LTW.Range("H140").Value = 87999.336
$ws_LTW.Range("J140").Value = 87999.5
$ws_LTW.Range("L140").Value = 87999.5
$ws_LTW.Range("N140").Value = -98359.5
$ws_WVR.Range("H18").Value = 507500
$ws_WVR.Range("I18").Value = 2000000
$ws_WVR.Range("J18").Value = 10000
$ws_WVR.Range("K18").Value = 2000000
$ws_WVR.Range("L18").Value = 10000
$ws_WVR.Range("M18").Value = -1999827
$ws_WVR.Range("N18").Value = -10346
$ws_WVR.Range("H24").Value = 11333.333
$ws_WVR.Range("I24").Value = 10000
$ws_WVR.Range("J24").Value = 12000
$ws_WVR.Range("K24").Value = 10000
$ws_WVR.Range("L24").Value = 12000
$ws_WVR.Range("N24").Value = -12460
$ws_WVR.Range("M24").Value = -9770
